$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 508, shifting existing rows 508:534 down to 509:535
$ws.Rows.Item(508).EntireRow.Insert()

# Populate the newly inserted row 508 with the new data record
$ws.Cells.Item(508, 1).Value = 3
$ws.Cells.Item(508, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(508, 3).Value = "Coquimbo"
$ws.Cells.Item(508, 4).Value = 45041
$ws.Cells.Item(508, 5).Value = 5
$ws.Cells.Item(508, 6).Value = 100112009
$ws.Cells.Item(508, 7).Value = "Acelga"
$ws.Cells.Item(508, 8).Value = "Sin especificar"
$ws.Cells.Item(508, 9).Value = "Primera"
$ws.Cells.Item(508, 10).Value = 240
$ws.Cells.Item(508, 11).Value = 3500
$ws.Cells.Item(508, 12).Value = 3800
$ws.Cells.Item(508, 13).Value = 3650
$ws.Cells.Item(508, 14).Value = "`$/docena de atados (6 kilos)"
$ws.Cells.Item(508, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(508, 16).Value = 608
$ws.Cells.Item(508, 17).Value = 6
$ws.Cells.Item(508, 18).Value = "Hortaliza"
